# Generate Report for Handback
# For the "14f697d3-bdef-4234-8c48-7bdad4f72596" row (row 6) on both the
# zh-cn and de-de sheets, record that a handback target/file was produced,
# but that its version is stale compared to the latest handoff.

$wb = $excel.ActiveWorkbook

$latestTargetDisplay = "14f697d3-bdef-4234-8c48-7bdad4f72596.md"
$latestTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77685ad8365c1c573236a5ac1cd668bf0af671c8/e2e/14f697d3-bdef-4234-8c48-7bdad4f72596.md"

function Update-HandbackRow {
    param(
        $ws,
        [string]$handbackFile,
        [string]$handbackDateTime
    )

    # I6: Latest Target File -> link to the latest handoff markdown file
    $ws.Range("I6").Value = $latestTargetDisplay
    $ws.Hyperlinks.Add($ws.Range("I6"), $latestTargetUrl, "", "", $latestTargetDisplay)
    $ws.Range("I6").Font.Underline = 2
    $ws.Range("I6").Font.Color = 15570276

    # J6: Latest Handback File -> the handed-back xlf (same as the handoff file, G6)
    $ws.Range("J6").Value = $handbackFile

    # K6: Latest Handback DateTime
    $ws.Range("K6").Value = $handbackDateTime

    # P6: Error Detail -> version mismatch message
    $ws.Range("P6").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc3c5228ed835efb96ed1611d5ae2ab987dbeb07/e2e/14f697d3-bdef-4234-8c48-7bdad4f72596.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77685ad8365c1c573236a5ac1cd668bf0af671c8/e2e/14f697d3-bdef-4234-8c48-7bdad4f72596.md."

    # Error Detail column needs to be wide enough to show the message.
    $ws.Columns.Item(16).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
}

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow -ws $wsZhCn -handbackFile "14f697d3-bdef-4234-8c48-7bdad4f72596.993963fd8d3edd7777ee6c5c965ff87f585acda6.zh-cn.xlf" -handbackDateTime "2016-09-01 16:51:45"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow -ws $wsDeDe -handbackFile "14f697d3-bdef-4234-8c48-7bdad4f72596.993963fd8d3edd7777ee6c5c965ff87f585acda6.de-de.xlf" -handbackDateTime "2016-09-01 16:51:52"
